{"js": "// Change the city/state shown in the resume's sidebar \"Text Box 2\" from\n// \"Kalamazoo, MI\" to \"Broomfield, CO\" (see commit message:\n// \"Update Resume to say Broomfield, CO instead of Kzoo\").\n//\n// The address line lives inside a floating text box (a DrawingML shape\n// anchored to the first heading paragraph), so it is not part of\n// context.document.body's own text -- it has to be reached through the\n// shape's own body/paragraphs.\n\nconst shapes = context.document.body.shapes;\nshapes.load(\"items/name\");\nawait context.sync();\n\n// Find the sidebar text box by name (falls back to the first shape that\n// contains the text we're looking for, in case the name ever changes).\nlet targetShape = null;\nfor (let i = 0; i < shapes.items.length; i++) {\n  if (shapes.items[i].name === \"Text Box 2\") {\n    targetShape = shapes.items[i];\n    break;\n  }\n}\nif (!targetShape && shapes.items.length > 0) {\n  targetShape = shapes.items[0];\n}\n\nif (targetShape) {\n  const shapeBody = targetShape.body;\n  const paragraphs = shapeBody.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (para.text.indexOf(\"Kalamazoo\") !== -1) {\n      const newText = para.text.replace(\"Kalamazoo, MI\", \"Broomfield, CO\");\n      para.insertText(newText, \"Replace\");\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Change the city/state shown in the resume's sidebar text box from\n# \"Kalamazoo, MI\" to \"Broomfield, CO\" (see commit message:\n# \"Update Resume to say Broomfield, CO instead of Kzoo\").\n#\n# The address line lives inside a floating text box (a drawing Shape\n# anchored to the first heading paragraph), so it has to be reached\n# through $d.Shapes rather than $d.Content.\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Shapes.Count; $i++) {\n    $shp = $d.Shapes.Item($i)\n    if ($shp.Name -eq \"Text Box 2\") {\n        $target = $shp\n        break\n    }\n}\nif ($target -eq $null -and $d.Shapes.Count -gt 0) {\n    $target = $d.Shapes.Item(1)\n}\n\nif ($target -ne $null) {\n    $rng = $target.TextFrame.TextRange\n\n    # Inspect the text through a duplicate range so the original $rng\n    # keeps pointing at just the address paragraph (first paragraph of the\n    # text box); only that duplicate's .Text needs to read the whole story.\n    $fullText = $rng.Duplicate.Text\n\n    if ($fullText.Contains(\"Kalamazoo, MI\")) {\n        $rng.Text = \"Broomfield, CO\"\n    }\n}\n"}
